# Insert a new data row at row 196 (pushes existing rows 196:324 down to
# 197:325, growing the used range from A1:R324 to A1:R325), then populate
# the newly-inserted row with a fresh "Apio" (celery) price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 196..324 down one row, inheriting formatting (e.g. the date
# number format on column D) from the row being pushed down.
$ws.Rows.Item(196).Insert()

$row = 196
$ws.Cells.Item($row, 1).Value  = 10
$ws.Cells.Item($row, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item($row, 3).Value  = "La Araucanía"
$ws.Cells.Item($row, 4).Value  = 44762
$ws.Cells.Item($row, 5).Value  = 9
$ws.Cells.Item($row, 6).Value  = 100112017
$ws.Cells.Item($row, 7).Value  = "Apio"
$ws.Cells.Item($row, 8).Value  = "Americana (o)"
$ws.Cells.Item($row, 9).Value  = "Primera"
$ws.Cells.Item($row, 10).Value = 45
$ws.Cells.Item($row, 11).Value = 11250
$ws.Cells.Item($row, 12).Value = 11250
$ws.Cells.Item($row, 13).Value = 11250
$ws.Cells.Item($row, 14).Value = "$/docena de matas"
$ws.Cells.Item($row, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 16).Value = 1875
$ws.Cells.Item($row, 17).Value = 6
$ws.Cells.Item($row, 18).Value = "Hortaliza"
